# Elaboración del Diagrama del proceso de manufactura manual de las bicicletas
# Adds an "Operario" (operator) column (F) and a workload-percentage column (G)
# to the manufacturing standard-time sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header cell F3 ("Operario") - reuse the same header style as C3/D3/E3
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("F3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F3").Value = "Operario"

# ---------------------------------------------------------------------------
# Operator letter codes for column F (rows 4-23) - reuse the bordered /
# centered body style already used by column B (style index 2)
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("F4:F23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Ordered the same way the original author entered them (non-sequential,
# matching the shared-string insertion order from the source edit) so the
# workbook's string table comes out byte-for-byte equivalent.
$fillOrder = @(
    @{ Row = 4;  Text = "A" },
    @{ Row = 20; Text = "O" },
    @{ Row = 5;  Text = "B" },
    @{ Row = 21; Text = "P" },
    @{ Row = 19; Text = "N" },
    @{ Row = 6;  Text = "C" },
    @{ Row = 17; Text = "L" },
    @{ Row = 18; Text = "M" },
    @{ Row = 7;  Text = "D, E" },
    @{ Row = 8;  Text = "F, G" },
    @{ Row = 9;  Text = "A" },
    @{ Row = 10; Text = "H" },
    @{ Row = 11; Text = "I" },
    @{ Row = 12; Text = "J" },
    @{ Row = 13; Text = "J" },
    @{ Row = 14; Text = "J" },
    @{ Row = 15; Text = "J" },
    @{ Row = 16; Text = "K" },
    @{ Row = 22; Text = "P" },
    @{ Row = 23; Text = "P" }
)

foreach ($entry in $fillOrder) {
    $ws.Range("F$($entry.Row)").Value = $entry.Text
}

# ---------------------------------------------------------------------------
# Workload percentage formulas for column G (rows 4-23)
# Rows 7 and 8 represent two operators sharing the task, so the load is
# halved; every other row is the straight percentage of an 8h (480 min) day.
# ---------------------------------------------------------------------------
$halfRows = @(7, 8)
for ($row = 4; $row -le 23; $row++) {
    $cell = $ws.Range("G$row")
    if ($halfRows -contains $row) {
        $cell.Formula = "=((D$row/480)*100)/2"
    } else {
        $cell.Formula = "=(D$row/480)*100"
    }
    $cell.NumberFormat = "0.00"
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
}

# ---------------------------------------------------------------------------
# Column G width (auto-fit-like sizing for the new percentage column)
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 8.8

# ---------------------------------------------------------------------------
# View state: zoom in and move the selection, matching the author's last
# on-screen position before saving.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 160
$ws.Range("G3").Select()
